# 30 Jan Presenti Sheet
# Add a new "30-Jan-24" attendance column (AG) to the Jan-2024 sheet,
# mirroring the formatting of the previous day's column (AF) and marking
# every student "Absent" for that day.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Jan-2024")

# New date header for 30-Jan-2024 in column AG, row 1 (same format as AF1)
$ws.Range("AF1").Copy()
$ws.Range("AG1").PasteSpecial(-4122)
$ws.Range("AG1").Value = [DateTime]::ParseExact("2024-01-30", "yyyy-MM-dd", $null)

# New attendance values for rows 2-4 in column AG (same format as AF2:AF4)
$ws.Range("AF2:AF4").Copy()
$ws.Range("AG2:AG4").PasteSpecial(-4122)
$ws.Range("AG2").Value = "Absent"
$ws.Range("AG3").Value = "Absent"
$ws.Range("AG4").Value = "Absent"

# Extend the data validation list to cover the new column
$ws.Range("C2:AF4").Validation.Delete()
$ws.Range("C2:AG4").Validation.Add(3, 1, 1, '"Present, Absent,Reason"')

# Update the current selection to reflect scrolling over to the new column
$ws.Range("AH7").Select()
